# Update countries & provincias Spain
# Refresh the COVID country table in the "Pais" sheet with the latest
# figures, and re-rank three countries (Polonia, Colombia, Niger) whose
# totals overtook the countries previously listed above them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($rowNumber, $values) {
    for ($i = 0; $i -lt $values.Count; $i++) {
        $ws.Cells.Item($rowNumber, $i + 1).Value = $values[$i]
    }
}

# --- Timestamp in the title cell (A1) -------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 3 de Abril de 2020 a las 23:20"

# --- Straight number refreshes (no re-ranking) -----------------------------
Set-Row 4   @("Estados Unidos", 272925, 28048, 12044, 253877, 5787, 934, 7004)
Set-Row 6   @("España",         119199,  7134, 30513,  77488, 6416, 850, 11198)
Set-Row 43  @("Grecia",           1613,    69,    78,   1472,   92,  10,    63)
Set-Row 104 @("Estado de Palestina", 194,   33,    21,    172,    0,   0,     1)
Set-Row 110 @("Georgia",            155,   21,    28,    127,    6,   0,     0)
Set-Row 147 @("Bermudas",            35,    0,    14,     21,    0,   0,     0)

# --- Polonia overtakes Ecuador & Malasia -----------------------------------
# New data puts Polonia ahead, pushing Ecuador and Malasia down one row.
Set-Row 32 @("Malasia", 3333, 217, 827, 2453, 108,  3,  53)
Set-Row 31 @("Ecuador", 3368, 205,  65, 3158, 100, 25, 145)
Set-Row 30 @("Polonia", 3383, 437,  56, 3256,  50, 14,  71)

# --- Colombia overtakes Argentina, Emiratos Arabes Unidos & Argelia --------
Set-Row 54 @("Argelia",                1171, 185,  62, 1004,   0, 19, 105)
Set-Row 53 @("Emiratos Arabes Unidos", 1264, 240, 108, 1147,   2,  1,   9)
Set-Row 52 @("Argentina",              1265,   0, 266,  960,   0,  3,  39)
Set-Row 51 @("Colombia",               1267, 106,  55, 1187,  50,  6,  25)

# --- Niger overtakes Isla de Man & Camboya ----------------------------------
Set-Row 122 @("Camboya",     114,  4, 35, 79, 1, 0, 0)
Set-Row 121 @("Isla de Man", 114, 19,  0, 113, 0, 0, 1)
Set-Row 120 @("Niger",       120, 22,  0, 115, 0, 0, 5)
